$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# This slide becomes a hidden slide in the slide show (adds show="0" on
# the slide's <p:sld> element).
$s.SlideShowTransition.Hidden = $true

# Locate the two "Subtitle 2" text boxes that need new copy by matching
# their current text, rather than relying on a fixed shape index.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $txt = $shp.TextFrame.TextRange.Text

        if ($txt -like "*Notify Node Agents about upcoming infrastructure updates*") {
            # New copy mentions applications too, so the box grows by one
            # line; set the autofit height PowerPoint would compute.
            $shp.TextFrame.TextRange.Text = "Notify Node Agents and applications about upcoming infrastructure updates  "
            $shp.Height = 61.96189
        }
        elseif ($txt -like "*Let Node Agents apply the scheduled changes*") {
            $shp.TextFrame.TextRange.Text = "Let Node Agents and application apply the scheduled changes at pre-determined times and collect reports"
        }
    }
}
